# Edit script: "Visions of the Quantum Realm" -> "The Enchanting Symphony of Colors"
# Applies the title/author/email changes, rewrites the body paragraph with the new
# "colors" essay content (incl. the Body:/Paragraph N: outline labels), rewrites the
# Summary paragraph, and appends a trailing empty paragraph.

$d = $word.ActiveDocument

# --- Title ---
$d.Paragraphs(1).Range.Text = "The Enchanting Symphony of Colors"

# --- Author name ---
$d.Paragraphs(2).Range.Text = "Roselyn Carter"

# --- Author email (kept as one paragraph: "user" + "." + "domain") ---
$d.Paragraphs(3).Range.Text = "roselyncarter@xyzschool.edu"

# --- Body paragraph (the long one with the line breaks) ---
$BR = [char]11
$bodyText = (
    "- Dive into the vibrant realm of colors, where hues dance together in a captivating symphony." +
    " From the blazing scarlet sunsets that ignite the skies to the calming cerulean depths of the ocean, colors enchant our world with their unspoken stories." +
    $BR + $BR +
    "- Colors, like musical notes, possess the power to evoke emotions and shape perceptions." +
    " From the vibrant reds that spark passion to the soothing greens that promote tranquility, colors have a profound impact on our psychological and physiological well-being." +
    $BR + $BR +
    "- Beyond their aesthetic allure, colors play a crucial role in various scientific fields." +
    " Whether it's the study of light and its interactions or the analysis of chemical compounds, colors serve as essential tools for unraveling the mysteries of the universe." +
    $BR + $BR +
    "Body:" +
    $BR + $BR +
    "Paragraph 1:" +
    $BR + $BR +
    "- In the realm of art, colors become the language of expression." +
    " Artists use colors to convey emotions, tell stories, and create visual masterpieces." +
    " From the bold strokes of abstract paintings to the intricate details of realistic landscapes, colors allow artists to share their unique perspectives and connect with viewers on a profound level." +
    $BR + $BR +
    "Paragraph 2:" +
    $BR + $BR +
    "- In the realm of science, colors hold significant importance." +
    " From the rainbow's spectrum, scientists gain insights into the properties of light and its interactions with matter." +
    " Colors also serve as indicators in chemical reactions, revealing the composition and structure of substances." +
    " Moreover, colors play a vital role in biotechnology and medical research, assisting in the development of diagnostic techniques and treatments." +
    $BR + $BR +
    "Paragraph 3:" +
    $BR + $BR +
    "- In the realm of history and culture, colors carry immense significance." +
    " Different cultures associate specific colors with emotions, values, and beliefs." +
    " From the vibrant colors of traditional festivals to the symbolic use of colors in flags and emblems, colors serve as cultural markers that reflect the heritage and identity of nations." +
    " Colors also play a crucial role in politics, influencing public perception and electoral outcomes."
)
$d.Paragraphs(5).Range.Text = $bodyText

# --- Summary paragraph ---
$summaryText = (
    "- Colors, with their inherent beauty and symbolism, hold a profound impact on our lives." +
    " From their role in art and science to their significance in culture and history, colors shape our perceptions, evoke emotions, and contribute to our understanding of the world around us." +
    " By delving into the fascinating world of colors, we unlock a treasure trove of knowledge, beauty, and inspiration."
)
$d.Paragraphs(7).Range.Text = $summaryText

# --- Trailing empty paragraph after the Summary body ---
$d.Paragraphs(7).Range.InsertParagraphAfter()
